$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.097.55'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +5.19%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.920.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.56%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -1.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.52%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9979'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5144'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3991'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.54%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08453'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '43.01'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.71%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.120'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.301'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.68%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.48%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.911.33'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.370'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9981'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.33%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.77%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001114'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06734'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.95'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.63%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9972'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.32%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.031'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.94%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.100.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.16%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.203'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.131.97'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.99%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.71%  '

$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.40'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.465'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.30%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.63'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.079'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.86%  '

$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.090'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.676'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.60%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02497'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06605'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.22%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.247'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.94%  '

$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2208'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.002'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.53%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.174'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.70%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6537'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.238'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.58%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.37'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6125'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.60%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.716'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.058'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.89%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.68%  '

$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.236'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.161'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.74%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '78.45'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.62%  '
